$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CSV block: A,C,D,E,F,G,H  for rows 2..31 (30 data rows).
# Column B is constant "walkingToRunning" for every row.
$csv = @"
0,-7.333617687225342,-4.549784660339356,0.0833368301391601,0.4665061936658925,-0.1984146457969052,0.8845337243641089
100,-15.08706188201904,-8.937231063842773,-1.605889558792114,0.7887173727446912,-0.2283256008636717,-2.244795449808533
200,11.69874000549316,-22.40719985961914,1.197917938232422,2.415960978059204,-0.9742195974378011,-3.041498450671924
300,3.118687152862549,-6.122324466705322,8.912668228149414,-5.063323306102351,-1.817073576590591,2.776077873566551
400,-8.666399955749512,-21.93803024291992,20.0407600402832,-6.820184679592355,0.3425412797460279,5.666934162962662
500,-24.77190589904785,16.44706535339355,-12.57914733886719,-0.2362654770122372,5.465962034814462,3.005506543552197
600,1.322300910949707,-0.346986711025238,-9.815587043762209,4.696316433887858,-2.277570799285295,-2.875306896134922
700,27.22784423828125,-27.38262748718262,-23.78179931640625,6.539746003992415,-0.2706960135815262,-1.030937138725735
800,49.64927673339844,-42.73597717285156,-6.524261474609375,4.692681789398195,1.697400997666751,-1.396730149493496
900,-3.964775562286377,4.51526403427124,1.280887007713318,-4.347883439531461,-0.5452939552419174,0.9300998636320361
1000,39.84653091430664,-19.79646873474121,38.97504425048828,-9.689171734978153,-3.703270496106596,9.968994430467122
1100,-19.31760787963867,17.29528427124023,-8.18376350402832,2.393433683058804,-9.479072037865135,1.899762265822456
1200,-16.73611831665039,-1.23279333114624,-18.61682319641113,4.808085488338101,-6.932880869098746,-1.150131797089299
1300,-40.93106842041016,-28.66251182556152,5.111406803131104,6.954113165537517,4.074747403462728,-1.481094896793366
1400,-5.353846549987793,3.492556095123291,14.51318454742432,0.2588006748873219,3.881287055857073,1.234817659153635
1500,-5.202733993530273,10.5782413482666,3.768710374832153,-3.794551531473794,-1.161946320066274,2.955449520372877
1600,17.96537399291992,-41.82672500610352,13.37196350097656,-3.552903596092673,6.386371462952866,1.291433391325612
1700,16.90632247924805,11.39815711975098,-9.546463966369627,-2.714477875653444,9.294961256139448,-3.878858893233145
1800,-13.79934120178223,-5.43239688873291,-14.25337219238281,2.263628258424627,-0.6571153846441176,-3.716604506268235
1900,-7.913525104522705,-43.03920745849609,5.675649166107178,3.927821926042158,-5.819554076475262,-0.2955599602531069
2000,-19.80719947814941,17.96865081787109,0.9991121292114258,2.031181216239927,-3.260402391938594,0.723365492680494
2100,5.307936668395996,-36.55269622802734,44.45186233520508,-5.857060266476065,-0.9245907862981151,6.962034623996895
2200,-72.68199157714844,-20.47943878173828,-39.62965774536133,-3.386064576167785,-0.2334743749861628,3.267812457739236
2300,-8.546328544616699,0.6573230028152466,-10.66314029693604,3.026514137492571,-1.510893523693044,-3.966962477740206
2400,-15.78201103210449,-11.37719249725342,-7.033545017242432,3.048832594179637,-2.553463636660109,-2.418893182978909
2500,-11.98593044281006,2.52801513671875,9.990962028503418,2.535124498255125,-1.166468461500665,-2.345369084208624
2600,-9.182360649108888,19.98668098449707,6.591425895690918,-3.627002197153437,-0.152242738972691,0.2601947556523713
2700,14.62302780151367,-34.03761672973633,14.29042530059814,-3.252175232943356,1.719280733197323,0.07690661035332244
2800,-34.80189895629883,14.59354972839356,-21.27434158325196,-2.007904552945881,8.391763182247416,-2.262497771020045
2900,-21.16343307495117,-7.425243377685547,-10.02677822113037,2.766276373582728,6.427707840414632,-4.763592523687029
"@

$lines = $csv -split "`n" | Where-Object { $_.Trim().Length -gt 0 }
$rowCount = $lines.Length

# Build 2D arrays (1-based indexing as used by Excel Range.Value arrays):
#  - $dataA  : single column, timestamps
#  - $dataB  : single column, constant label
#  - $dataCH : six columns, ax,ay,az,gx,gy,gz
$dataA = New-Object 'object[,]' $rowCount,1
$dataB = New-Object 'object[,]' $rowCount,1
$dataCH = New-Object 'object[,]' $rowCount,6

for ($i = 0; $i -lt $rowCount; $i++) {
    $parts = $lines[$i].Trim() -split ","
    $dataA[$i, 0] = [double]$parts[0]
    $dataB[$i, 0] = "walkingToRunning"
    for ($j = 0; $j -lt 6; $j++) {
        $dataCH[$i, $j] = [double]$parts[$j + 1]
    }
}

$lastRow = 1 + $rowCount   # header is row 1, data starts row 2

$ws.Range("A2:A$lastRow").Value = $dataA
$ws.Range("B2:B$lastRow").Value = $dataB
$ws.Range("C2:H$lastRow").Value = $dataCH
